$d = $word.ActiveDocument

# --- Fix footnote 28's text (was just the danda "།", becomes the full note) ---
# Footnotes collection is 1-based by position; footnote w:id="28" is the 8th
# footnote in this document (ids 21..35 map to Index 1..15).
$fn28 = $d.Footnotes.Item(8)
$fn28.Range.Text = "ལ། ཞེས་པར་མ་གཞན་ནང་མེད།"

# --- Fix footnote 34's text: drop the stray trailing "aa" typo ---
$fn34 = $d.Footnotes.Item(14)
$fn34.Range.Text = "ཀླུ་གྲུབ། སྣར་ཐང་།"

# --- Remove the empty footnote 35 entirely (the "empty notes" bug) ---
# Deleting the Footnote object removes both its footnoteReference run in the
# body and its <w:footnote> entry, renumbering the remaining ones.
$fn35 = $d.Footnotes.Item(15)
$fn35.Delete()
